# Generate Report for Handoff
#
# Stamps the "Latest Handoff Datetime" for every file that is currently
# queued for handoff (status "Ready for handoff") or needs to be
# re-handed-off because the previous handback transform failed
# ("Handback transform failed"), with the timestamp of this handoff run.
# Files that are still "In Translation" or have already been
# "Handed back: in sync with en-US" are left untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# New handoff timestamps generated by this run, one per locale sheet plus
# the rolled-up value shown on the Overview sheet.
$overviewStamp = "2016-24-19 02:24:37"
$zhcnStamp     = "2016-03-19 02:24:34"
$dedeStamp     = "2016-03-19 02:24:37"

# Rows on the per-locale sheets (and the matching Overview row) whose
# handoff datetime gets refreshed by this report run.
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = $overviewStamp
    $zhcn.Range("E$r").Value = $zhcnStamp
    $dede.Range("E$r").Value = $dedeStamp
}
